$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.556.85"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.629.13"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'212.70"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'0.250"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "1.855.81"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.646.00"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'4.12"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "'63.85"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "26.605.38"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "'214.81"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").Value = "'1.00"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "'4.31"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("E24").Value = "  +4.27%  "
$ws.Range("D25").Value = "'148.28"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'6.86"
$ws.Range("E28").Value = "  +3.36%  "
$ws.Range("D29").Value = "'15.52"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "'0.0506"
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "1.220.58"
$ws.Range("E35").Value = "  +4.96%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'0.0173"
$ws.Range("E37").Value = "  +5.11%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").Value = "'0.796"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("D42").Value = "'0.793"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "'5.36"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "1.767.04"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'92.30"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'54.91"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  +6.01%  "
$ws.Range("D49").Value = "'0.0511"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'7.64"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  -0.12%  "
